# Generate Report for Handoff
# Update status text and timestamps across the Overview, zh-cn, and de-de
# sheets, and widen the affected "datetime" columns to fit the new values.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet ---------------------------------------------------
# E2 (zh-cn status) and F2 (de-de status): "In Translation" -> "Ready for handoff"
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
# G2 (Latest HO Xliff Generate Date)
$wsOverview.Range("G2").Value = "2016-09-02 19:06:34"

# Widen columns E and F to fit the longer status text
# (ColumnWidth is specified in character units; 16.3333 is the input that
# yields the target stored column width used by the workbook.)
$wsOverview.Columns.Item(5).ColumnWidth = 16.333333333333332
$wsOverview.Columns.Item(6).ColumnWidth = 16.333333333333332

# --- zh-cn sheet --------------------------------------------------------
# C2 (Status): "In Translation" -> "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
# H2 (Latest Handoff Datetime)
$wsZhCn.Range("H2").Value = "2016-09-02 19:06:30"

# Widen column C to fit the longer status text
$wsZhCn.Columns.Item(3).ColumnWidth = 16.333333333333332

# --- de-de sheet --------------------------------------------------------
# C2 (Status): "In Translation" -> "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"
# H2 (Latest Handoff Datetime / Latest HO Xliff Generate Date equivalent)
$wsDeDe.Range("H2").Value = "2016-09-02 19:06:34"

# Widen column C to fit the longer status text
$wsDeDe.Columns.Item(3).ColumnWidth = 16.333333333333332
